$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(85).Insert(1)  # xlShiftDown = -4121, but default is down when inserting a row

$ws.Range("A85").Value = 11
$ws.Range("B85").Value = "Vega Monumental Concepción"
$ws.Range("C85").Value = "Bíobío"
$ws.Range("D85").Value = 44663
$ws.Range("E85").Value = 8
$ws.Range("F85").Value = "Fruta"
$ws.Range("G85").Value = 100101
$ws.Range("H85").Value = "Berries"
$ws.Range("I85").Value = 100101007
$ws.Range("J85").Value = "Kiwi"
$ws.Range("K85").Value = "Hayward"
$ws.Range("L85").Value = "Primera"
$ws.Range("M85").Value = 180
$ws.Range("N85").Value = 10000
$ws.Range("O85").Value = 11000
$ws.Range("P85").Value = 10444
$ws.Range("Q85").Value = "$/bandeja 18 kilos"
$ws.Range("R85").Value = "Provincia de Curicó"
$ws.Range("S85").Value = 580
$ws.Range("T85").Value = 18
